$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update rows 2-63 with new contract data
$ws.Cells.Item(2, 1).Value = 'IUniswapV3FlashCallback'
$ws.Cells.Item(2, 2).Value = 0
$ws.Cells.Item(2, 3).Value = 0
$ws.Cells.Item(3, 1).Value = 'IUniswapV3MintCallback'
$ws.Cells.Item(3, 2).Value = 0
$ws.Cells.Item(3, 3).Value = 0
$ws.Cells.Item(4, 1).Value = 'IUniswapV3SwapCallback'
$ws.Cells.Item(4, 2).Value = 0
$ws.Cells.Item(4, 3).Value = 0
$ws.Cells.Item(5, 1).Value = 'IERC20Minimal'
$ws.Cells.Item(5, 2).Value = 0
$ws.Cells.Item(5, 3).Value = 0
$ws.Cells.Item(6, 1).Value = 'IUniswapV3Factory'
$ws.Cells.Item(6, 2).Value = 0
$ws.Cells.Item(6, 3).Value = 0
$ws.Cells.Item(7, 1).Value = 'IUniswapV3Pool'
$ws.Cells.Item(7, 2).Value = 6
$ws.Cells.Item(7, 3).Value = 0
$ws.Cells.Item(8, 1).Value = 'IUniswapV3PoolDeployer'
$ws.Cells.Item(8, 2).Value = 0
$ws.Cells.Item(8, 3).Value = 0
$ws.Cells.Item(9, 1).Value = 'IUniswapV3PoolActions'
$ws.Cells.Item(9, 2).Value = 0
$ws.Cells.Item(9, 3).Value = 0
$ws.Cells.Item(10, 1).Value = 'IUniswapV3PoolDerivedState'
$ws.Cells.Item(10, 2).Value = 0
$ws.Cells.Item(10, 3).Value = 0
$ws.Cells.Item(11, 1).Value = 'IUniswapV3PoolEvents'
$ws.Cells.Item(11, 2).Value = 0
$ws.Cells.Item(11, 3).Value = 0
$ws.Cells.Item(12, 1).Value = 'IUniswapV3PoolImmutables'
$ws.Cells.Item(12, 2).Value = 0
$ws.Cells.Item(12, 3).Value = 0
$ws.Cells.Item(13, 1).Value = 'IUniswapV3PoolOwnerActions'
$ws.Cells.Item(13, 2).Value = 0
$ws.Cells.Item(13, 3).Value = 0
$ws.Cells.Item(14, 1).Value = 'IUniswapV3PoolState'
$ws.Cells.Item(14, 2).Value = 0
$ws.Cells.Item(14, 3).Value = 0
$ws.Cells.Item(15, 1).Value = 'BitMath'
$ws.Cells.Item(15, 2).Value = 0
$ws.Cells.Item(15, 3).Value = 0
$ws.Cells.Item(16, 1).Value = 'FixedPoint128'
$ws.Cells.Item(16, 2).Value = 0
$ws.Cells.Item(16, 3).Value = 0
$ws.Cells.Item(17, 1).Value = 'FixedPoint96'
$ws.Cells.Item(17, 2).Value = 0
$ws.Cells.Item(17, 3).Value = 0
$ws.Cells.Item(18, 1).Value = 'FullMath'
$ws.Cells.Item(18, 2).Value = 0
$ws.Cells.Item(18, 3).Value = 0
$ws.Cells.Item(19, 1).Value = 'LiquidityMath'
$ws.Cells.Item(19, 2).Value = 0
$ws.Cells.Item(19, 3).Value = 0
$ws.Cells.Item(20, 1).Value = 'LowGasSafeMath'
$ws.Cells.Item(20, 2).Value = 0
$ws.Cells.Item(20, 3).Value = 0
$ws.Cells.Item(21, 1).Value = 'Oracle'
$ws.Cells.Item(21, 2).Value = 0
$ws.Cells.Item(21, 3).Value = 0
$ws.Cells.Item(22, 1).Value = 'Position'
$ws.Cells.Item(22, 2).Value = 3
$ws.Cells.Item(22, 3).Value = 0
$ws.Cells.Item(23, 1).Value = 'SafeCast'
$ws.Cells.Item(23, 2).Value = 0
$ws.Cells.Item(23, 3).Value = 0
$ws.Cells.Item(24, 1).Value = 'SqrtPriceMath'
$ws.Cells.Item(24, 2).Value = 5
$ws.Cells.Item(24, 3).Value = 0
$ws.Cells.Item(25, 1).Value = 'SwapMath'
$ws.Cells.Item(25, 2).Value = 2
$ws.Cells.Item(25, 3).Value = 0
$ws.Cells.Item(26, 1).Value = 'Tick'
$ws.Cells.Item(26, 2).Value = 4
$ws.Cells.Item(26, 3).Value = 0
$ws.Cells.Item(27, 1).Value = 'TickBitmap'
$ws.Cells.Item(27, 2).Value = 1
$ws.Cells.Item(27, 3).Value = 0
$ws.Cells.Item(28, 1).Value = 'TickMath'
$ws.Cells.Item(28, 2).Value = 0
$ws.Cells.Item(28, 3).Value = 0
$ws.Cells.Item(29, 1).Value = 'TransferHelper'
$ws.Cells.Item(29, 2).Value = 1
$ws.Cells.Item(29, 3).Value = 0
$ws.Cells.Item(30, 1).Value = 'UnsafeMath'
$ws.Cells.Item(30, 2).Value = 0
$ws.Cells.Item(30, 3).Value = 0
$ws.Cells.Item(31, 1).Value = 'NoDelegateCall'
$ws.Cells.Item(31, 2).Value = 0
$ws.Cells.Item(31, 3).Value = 0
$ws.Cells.Item(32, 1).Value = 'BitMathEchidnaTest'
$ws.Cells.Item(32, 2).Value = 1
$ws.Cells.Item(32, 3).Value = 0
$ws.Cells.Item(33, 1).Value = 'BitMathTest'
$ws.Cells.Item(33, 2).Value = 1
$ws.Cells.Item(33, 3).Value = 0
$ws.Cells.Item(34, 1).Value = 'FullMathEchidnaTest'
$ws.Cells.Item(34, 2).Value = 1
$ws.Cells.Item(34, 3).Value = 0
$ws.Cells.Item(35, 1).Value = 'FullMathTest'
$ws.Cells.Item(35, 2).Value = 1
$ws.Cells.Item(35, 3).Value = 0
$ws.Cells.Item(36, 1).Value = 'LiquidityMathTest'
$ws.Cells.Item(36, 2).Value = 1
$ws.Cells.Item(36, 3).Value = 0
$ws.Cells.Item(37, 1).Value = 'LowGasSafeMathEchidnaTest'
$ws.Cells.Item(37, 2).Value = 1
$ws.Cells.Item(37, 3).Value = 0
$ws.Cells.Item(38, 1).Value = 'MockTimeUniswapV3Pool'
$ws.Cells.Item(38, 2).Value = 1
$ws.Cells.Item(38, 3).Value = 0
$ws.Cells.Item(39, 1).Value = 'MockTimeUniswapV3PoolDeployer'
$ws.Cells.Item(39, 2).Value = 2
$ws.Cells.Item(39, 3).Value = 0
$ws.Cells.Item(40, 1).Value = 'NoDelegateCallTest'
$ws.Cells.Item(40, 2).Value = 1
$ws.Cells.Item(40, 3).Value = 0
$ws.Cells.Item(41, 1).Value = 'OracleEchidnaTest'
$ws.Cells.Item(41, 2).Value = 1
$ws.Cells.Item(41, 3).Value = 0
$ws.Cells.Item(42, 1).Value = 'OracleTest'
$ws.Cells.Item(42, 2).Value = 1
$ws.Cells.Item(42, 3).Value = 0
$ws.Cells.Item(43, 1).Value = 'SqrtPriceMathEchidnaTest'
$ws.Cells.Item(43, 2).Value = 3
$ws.Cells.Item(43, 3).Value = 0
$ws.Cells.Item(44, 1).Value = 'SqrtPriceMathTest'
$ws.Cells.Item(44, 2).Value = 1
$ws.Cells.Item(44, 3).Value = 0
$ws.Cells.Item(45, 1).Value = 'SwapMathEchidnaTest'
$ws.Cells.Item(45, 2).Value = 1
$ws.Cells.Item(45, 3).Value = 0
$ws.Cells.Item(46, 1).Value = 'SwapMathTest'
$ws.Cells.Item(46, 2).Value = 1
$ws.Cells.Item(46, 3).Value = 0
$ws.Cells.Item(47, 1).Value = 'TestERC20'
$ws.Cells.Item(47, 2).Value = 1
$ws.Cells.Item(47, 3).Value = 0
$ws.Cells.Item(48, 1).Value = 'TestUniswapV3Callee'
$ws.Cells.Item(48, 2).Value = 7
$ws.Cells.Item(48, 3).Value = 0
$ws.Cells.Item(49, 1).Value = 'TestUniswapV3ReentrantCallee'
$ws.Cells.Item(49, 2).Value = 3
$ws.Cells.Item(49, 3).Value = 0
$ws.Cells.Item(50, 1).Value = 'TestUniswapV3Router'
$ws.Cells.Item(50, 2).Value = 5
$ws.Cells.Item(50, 3).Value = 0
$ws.Cells.Item(51, 1).Value = 'TestUniswapV3SwapPay'
$ws.Cells.Item(51, 2).Value = 3
$ws.Cells.Item(51, 3).Value = 0
$ws.Cells.Item(52, 1).Value = 'TickBitmapEchidnaTest'
$ws.Cells.Item(52, 2).Value = 1
$ws.Cells.Item(52, 3).Value = 0
$ws.Cells.Item(53, 1).Value = 'TickBitmapTest'
$ws.Cells.Item(53, 2).Value = 1
$ws.Cells.Item(53, 3).Value = 0
$ws.Cells.Item(54, 1).Value = 'TickEchidnaTest'
$ws.Cells.Item(54, 2).Value = 1
$ws.Cells.Item(54, 3).Value = 0
$ws.Cells.Item(55, 1).Value = 'TickMathEchidnaTest'
$ws.Cells.Item(55, 2).Value = 1
$ws.Cells.Item(55, 3).Value = 0
$ws.Cells.Item(56, 1).Value = 'TickMathTest'
$ws.Cells.Item(56, 2).Value = 1
$ws.Cells.Item(56, 3).Value = 0
$ws.Cells.Item(57, 1).Value = 'TickOverflowSafetyEchidnaTest'
$ws.Cells.Item(57, 2).Value = 1
$ws.Cells.Item(57, 3).Value = 0
$ws.Cells.Item(58, 1).Value = 'TickTest'
$ws.Cells.Item(58, 2).Value = 1
$ws.Cells.Item(58, 3).Value = 0
$ws.Cells.Item(59, 1).Value = 'UniswapV3PoolSwapTest'
$ws.Cells.Item(59, 2).Value = 3
$ws.Cells.Item(59, 3).Value = 0
$ws.Cells.Item(60, 1).Value = 'UnsafeMathEchidnaTest'
$ws.Cells.Item(60, 2).Value = 1
$ws.Cells.Item(60, 3).Value = 0
$ws.Cells.Item(61, 1).Value = 'UniswapV3Factory'
$ws.Cells.Item(61, 2).Value = 4
$ws.Cells.Item(61, 3).Value = 0
$ws.Cells.Item(62, 1).Value = 'UniswapV3Pool'
$ws.Cells.Item(62, 2).Value = 21
$ws.Cells.Item(62, 3).Value = 0
$ws.Cells.Item(63, 1).Value = 'UniswapV3PoolDeployer'
$ws.Cells.Item(63, 2).Value = 2
$ws.Cells.Item(63, 3).Value = 0

# Delete rows 64-77 (14 rows that no longer exist)
$ws.Range("A64:C77").EntireRow.Delete()

Write-Host "Done. UsedRange rows: $($ws.UsedRange.Rows.Count())"
